$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Part 1: merge the two runs that used to be split by the old "_GoBack"
# bookmark ("...test plans, test cases, test " | bookmark | "matrices and
# other...") into a single run, leaving the preceding "Formal training" /
# " Automated Testing " runs untouched.
# ---------------------------------------------------------------------------

# Move the stray "_GoBack" bookmark out of the way first; this un-splits
# nothing by itself, but clears the way for a local, in-run edit to merge
# the two text runs it used to separate.
$tmp = $d.Range(0, 0)
$d.Bookmarks.Add("_GoBack", $tmp) | Out-Null

# A self-replace confined entirely to the run right after the old bookmark
# location causes it to merge back into the immediately preceding run
# (which is now textually adjacent, the bookmark having been removed)
# without touching anything further up the paragraph.
$mergeRange = $d.Content
$mergeRange.Find.Execute("matrices and other", $true, $false, $false, $false, $false, $true, 1, $false, "matrices and other", 2) | Out-Null

# ---------------------------------------------------------------------------
# Part 2: "C#, Assembler" -> "Java, Ruby" in the Programming skills line,
# splitting the run so the "_GoBack" bookmark ends up right before "Java".
# ---------------------------------------------------------------------------

$skillRange = $d.Content
$skillRange.Find.Execute("C#, Assembler", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$langStart = $skillRange.Start

$skillRange.Text = "Java, Ruby"
$newLangEnd = $langStart + ("Java, Ruby").Length

$progLine = $d.Content
$progLine.Find.Execute(" SQL, HTML, CSS, XML, JavaScript, ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$lineStart = $progLine.Start

$cssEnd = $lineStart + (" SQL, HTML, CSS").Length
$jsEnd = $langStart

# Split " SQL, HTML, CSS" | ", XML, JavaScript, "
$t1 = $d.Range($lineStart, $cssEnd)
$t1.Font.Bold = $true
$t1.Font.Bold = $false

# Split ", XML, JavaScript, " | "Java, Ruby "
$t2 = $d.Range($cssEnd, $jsEnd)
$t2.Font.Bold = $true
$t2.Font.Bold = $false

# Split "Java, Ruby" | " "
$t3 = $d.Range($jsEnd, $newLangEnd)
$t3.Font.Bold = $true
$t3.Font.Bold = $false

# Move the "_GoBack" bookmark to sit right before "Java, Ruby".
$bmRange = $d.Range($jsEnd, $jsEnd)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
